# Add a new "2021" column (column R) to the SDG 1.5.1 indicator sheet,
# mirroring the formatting of the existing 2020 column (Q) and filling in
# the reported values (missing data points use the existing "-" placeholder
# shared string already used elsewhere in the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of column Q (borders/fonts/number format/alignment)
# into the new column R for every row that currently has data in Q
# (the thin top border row 3 through the last data row 34).
$ws.Range("Q3:Q34").Copy()
$ws.Range("R3:R34").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Year header
$ws.Range("R4").Value = 2021

# Кыргыз Республикасы / Kyrgyz Republic totals
$ws.Range("R5").Value = 109
$ws.Range("R6").Value = 74
$ws.Range("R7").Value = 35

# Баткен облусу / Batken oblast
$ws.Range("R8").Value = 36
$ws.Range("R9").Value = 35
$ws.Range("R10").Value = 1

# Жалал-Абад облусу / Djalal-Abad oblast
$ws.Range("R11").Value = 15
$ws.Range("R12").Value = 8
$ws.Range("R13").Value = 7

# Ысык-Көл облусу / Ysyk-Kul oblast
$ws.Range("R14").Value = 12
$ws.Range("R15").Value = 7
$ws.Range("R16").Value = 5

# Нарын облусу / Naryn oblast - no data reported
$ws.Range("R17").Value = "-"
$ws.Range("R18").Value = "-"
$ws.Range("R19").Value = "-"

# Ош облусу / Osh oblast
$ws.Range("R20").Value = 17
$ws.Range("R21").Value = 8
$ws.Range("R22").Value = 9

# Талас облусу / Talas oblast
$ws.Range("R23").Value = 9
$ws.Range("R24").Value = 7
$ws.Range("R25").Value = 2

# Чүй облусу / Chui oblast
$ws.Range("R26").Value = 20
$ws.Range("R27").Value = 9
$ws.Range("R28").Value = 11

# Бишкек ш. / Bishkek city - no data reported
$ws.Range("R29").Value = "-"
$ws.Range("R30").Value = "-"
$ws.Range("R31").Value = "-"

# Ош ш. / Osh city - no data reported
$ws.Range("R32").Value = "-"
$ws.Range("R33").Value = "-"
$ws.Range("R34").Value = "-"

# Leave the cursor on A1, matching the saved file's lack of an explicit
# (non-default) active-cell selection.
$ws.Range("A1").Select() | Out-Null
